$wb = $excel.ActiveWorkbook

# --- 1. Fix ordering of the cfop list string on "PI hours" (G5) ---
$piWs = $wb.Worksheets.Item("PI hours")
$piWs.Range("G5").Value = "['cfop_RRC', 'cfop_CHOUDHURY']"

# --- 2. Add the "users" column (E) to "project hours" ---
$projWs = $wb.Worksheets.Item("project hours")

$projWs.Range("E1").Value = "users"
$projWs.Range("B1").Copy()
$projWs.Range("E1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$projWs.Range("E2").Value = "['Hunter Young', 'Harshal Maske']"
$projWs.Range("E3").Value = "['Usman syed', 'Usman Syed']"
$projWs.Range("E4").Value = "['Arun Lakshmanan', 'Mitchell Jones']"
$projWs.Range("E5").Value = "['Harshal Maske']"
$projWs.Range("E6").Value = "['Jonathan Hoff']"
$projWs.Range("E7").Value = "['Sheng Shen']"
